# Auto-generated edit script applying cell-level changes from the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.924.30'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '2.196.89'
$ws.Range('E3').Value = '  -2.31%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''295.22'
$ws.Range('E5').Value = '  -4.02%  '
$ws.Range('D6').Value = '''88.77'
$ws.Range('E6').Value = '  -6.45%  '
$ws.Range('D7').Value = '''0.566'
$ws.Range('E7').Value = '  -0.83%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '''0.482'
$ws.Range('E9').Value = '  -8.50%  '
$ws.Range('D10').Value = '''32.04'
$ws.Range('E10').Value = '  -8.11%  '
$ws.Range('D11').Value = '''0.0771'
$ws.Range('E11').Value = '  -4.96%  '
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('D13').Value = '''6.78'
$ws.Range('E13').Value = '  -6.07%  '
$ws.Range('D14').Value = '2.532.95'
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').Value = '2.269.90'
$ws.Range('E15').Value = '  -6.12%  '
$ws.Range('D16').Value = '''13.07'
$ws.Range('E16').Value = '  -4.77%  '
$ws.Range('D17').Value = '''0.772'
$ws.Range('E17').Value = '  -8.21%  '
$ws.Range('D18').Value = '43.595.59'
$ws.Range('E18').Value = '  -1.16%  '
$ws.Range('D19').Value = '0.0₃0890'
$ws.Range('E19').Value = '  -7.72%  '
$ws.Range('D20').Value = '''5.84'
$ws.Range('E20').Value = '  -8.86%  '
$ws.Range('D21').Value = '''10.80'
$ws.Range('E21').Value = '  -13.70%  '
$ws.Range('D22').Value = '''62.98'
$ws.Range('E22').Value = '  -4.49%  '
$ws.Range('D23').Value = '''231.05'
$ws.Range('E23').Value = '  -2.88%  '
$ws.Range('E24').Value = '  -8.21%  '
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('E26').Value = '  -9.10%  '
$ws.Range('D27').Value = '''2.23'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').Value = '''35.92'
$ws.Range('E28').Value = '  -6.52%  '
$ws.Range('D29').Value = '''9.22'
$ws.Range('E29').Value = '  -6.86%  '
$ws.Range('D30').Value = '''19.13'
$ws.Range('E30').Value = '  -4.88%  '
$ws.Range('D31').Value = '''147.05'
$ws.Range('E31').Value = '  -4.55%  '
$ws.Range('D32').Value = '''5.30'
$ws.Range('E32').Value = '  -11.27%  '
$ws.Range('D33').Value = '''2.51'
$ws.Range('E33').Value = '  -5.33%  '
$ws.Range('D34').Value = '''0.0735'
$ws.Range('E34').Value = '  -8.50%  '
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('D36').Value = '''2.88'
$ws.Range('E36').Value = '  -7.34%  '
$ws.Range('E37').Value = '  -6.49%  '
$ws.Range('D38').Value = '''1.65'
$ws.Range('E38').Value = '  -9.14%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '''3.10'
$ws.Range('E39').Value = '  -11.13%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.0282'
$ws.Range('E40').Value = '  -7.42%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '''3.52'
$ws.Range('E41').Value = '  -8.09%  '
$ws.Range('D42').Value = '''1.01'
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('D43').Value = '''13.07'
$ws.Range('E43').Value = '  -11.95%  '
$ws.Range('D44').Value = '1.792.02'
$ws.Range('E44').Value = '  +2.37%  '
$ws.Range('D45').Value = '''1.65'
$ws.Range('E45').Value = '  +1.76%  '
$ws.Range('E46').Value = '  +11.46%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '''0.173'
$ws.Range('E47').Value = '  -10.79%  '
$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').Value = '''72.54'
$ws.Range('E48').Value = '  -10.58%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '''91.83'
$ws.Range('E49').Value = '  -8.10%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').Value = '''65.23'
$ws.Range('E50').Value = '  -7.87%  '
$ws.Range('D51').Value = '2.415.38'
$ws.Range('E51').Value = '  -2.21%  '

# Reset style on cells we force-quoted as text back to Normal,
# so no stray 'quote prefix' cell style lingers (keeps formatting identical to source).
$textCells = @('D5','D6','D7','D9','D10','D11','D13','D16','D17','D20','D21','D22','D23','D27','D28','D29','D30','D31','D32','D33','D34','D36','D38','D39','D40','D41','D42','D43','D45','D47','D48','D49','D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = 'Normal'
}
